# ContactPage locator dictionary update
# Commit: "Automate remaining test cases" — the Privacy Consent checkbox/label
# locators moved from fixed Marketo-generated ids (input#mktoCheckbox_142098_0 /
# label#LblmktoCheckbox_142098_0) to stable structural selectors
# (div.mktoCheckboxList input / div.mktoCheckboxList label).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ContactPage")

# Row 27: ContactPage_CheckBox_PrivacyConsent locator value
$ws.Range("C27").Value = "div.mktoCheckboxList input"

# Row 28: ContactPage_Label_PrivacyConsent locator value
$ws.Range("C28").Value = "div.mktoCheckboxList label"

# Leave the sheet with the same selection/zoom state recorded for this edit
$ws.Range("C13").Select()
$excel.ActiveWindow.Zoom = 110
